$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update harvester column (B) for all data rows (2-25) to the new value,
# and fill in the experimentDesign column (D) with the new value.
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 2).Value = "S.GISH"
    $ws.Cells.Item($r, 4).Value = "90minuteInduction"
}

# Match the selection left behind in the saved workbook (D2:D25, active cell D2)
$ws.Range("D2:D25").Select()
